{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright block (and the blank\n// paragraph that preceded it) that used to follow the \"Requisitos\" list,\n// while keeping the paragraph(s) after it (blank paragraph + page-break\n// paragraph at the very end of the document body) intact.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOM3234: \u00d3ptica F\u00edsica (Requisito)\") and the\n// \"\u00a9 2020 ...\" paragraph that marks the end of the block to remove.\nlet anchorIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"LOM3234\") !== -1) {\n    anchorIndex = i;\n  }\n  if (text.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1 && copyrightIndex !== -1 && copyrightIndex > anchorIndex) {\n  // Delete every paragraph strictly between the anchor and (inclusive) the\n  // copyright paragraph, i.e. the blank line + \"Ver no Jupiter...\" line +\n  // \"\u00a9 2020 ...\" line.\n  for (let i = copyrightIndex; i > anchorIndex; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright block (and the blank\n# paragraph that preceded it) that used to follow the \"Requisitos\" list,\n# while keeping the paragraph(s) after it (blank paragraph + page-break\n# paragraph at the very end of the document body) intact.\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\n$copyrightIndex = -1\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($text -like \"*LOM3234*\") {\n        $anchorIndex = $i\n    }\n    if ($text -like \"*Contact: luizeleno@usp.br*\") {\n        $copyrightIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1 -and $copyrightIndex -ne -1 -and $copyrightIndex -gt $anchorIndex) {\n    # Delete every paragraph strictly between the anchor and (inclusive) the\n    # copyright paragraph, i.e. the blank line + \"Ver no Jupiter...\" line +\n    # \"\u00a9 2020 ...\" line. Walk backwards so indices stay valid as we delete.\n    for ($i = $copyrightIndex; $i -gt $anchorIndex; $i--) {\n        $d.Paragraphs.Item($i).Range.Delete()\n    }\n}\n"}
